# Fix "Recorded By" (column G) entries so that the literal "System" entry
# is no longer listed first in the comma-separated list of recorders.
# For every row where column G starts with "System, ", swap the first two
# comma-separated tokens (i.e. move "System" to the second position),
# leaving any remaining tokens untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G = 7
    $text = $cell.Value()

    if ($null -eq $text) {
        continue
    }

    $value = [string]$text

    if ($value.StartsWith("System, ")) {
        $parts = $value.Split(",")
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }

        if ($parts.Length -ge 2) {
            $first = $parts[0]
            $second = $parts[1]
            $parts[0] = $second
            $parts[1] = $first
            $newValue = [string]::Join(", ", $parts)
            $cell.Value = $newValue
        }
    }
}
